# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: EDWIN ENRIQUE RODRIGUEZ PEREZ - periodo 1803 (new row moved to top)
$ws.Range("C16").Value = "1077453261"
$ws.Range("D16").Value = "EDWIN ENRIQUE RODRIGUEZ PEREZ"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 33600
$ws.Range("G16").Value = 953370

# Row 17: EDWIN ENRIQUE RODRIGUEZ PEREZ - periodo 1804
$ws.Range("C17").Value = "1077453261"
$ws.Range("D17").Value = "EDWIN ENRIQUE RODRIGUEZ PEREZ"
$ws.Range("E17").Value = "1804"
$ws.Range("F17").Value = 33600
$ws.Range("G17").Value = 953370

# Row 18: ANDRES ALFONSO MACKENZIE LOPEZ - periodo 1808 (salario actualizado)
$ws.Range("C18").Value = "1066517619"
$ws.Range("D18").Value = "ANDRES ALFONSO MACKENZIE LOPEZ"
$ws.Range("E18").Value = "1808"
$ws.Range("F18").Value = 32000
$ws.Range("G18").Value = 2209505

# Row 19: MARITZA HOYOS BARRAZA - periodo 1810
$ws.Range("C19").Value = "22789256"
$ws.Range("D19").Value = "MARITZA HOYOS BARRAZA"
$ws.Range("E19").Value = "1810"
$ws.Range("F19").Value = 34000
$ws.Range("G19").Value = 964080
